$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to reflect the new "through" date
$ws.Name = "Through 2021-12-05"

# Update the December partial-month label
$ws.Range("A13").Value = "December (through 12-05)"

# November (row 12) - 2021 column (H) correction
$ws.Range("H12").Value = 201

# December (row 13) - new data for 2021-12-13
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = 16
$ws.Range("D13").Value = 19
$ws.Range("E13").Value = 11
$ws.Range("F13").Value = 5
$ws.Range("G13").Value = 26
$ws.Range("H13").Value = 37

# Total (row 14) - updated sums
$ws.Range("B14").Value = 294
$ws.Range("C14").Value = 579
$ws.Range("D14").Value = 840
$ws.Range("E14").Value = 693
$ws.Range("F14").Value = 539
$ws.Range("G14").Value = 1290
$ws.Range("H14").Value = 1680
